# Re-run SGNN to annotate dialog acts following clean up work to the original transcripts.
# Updates the DAMSLTag (column I) and DialogAct (column J) values for a number of rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> (DAMSLTag, DialogAct)
$updates = @{
    3   = @("b",  "Acknowledge (Backchannel)")
    4   = @("b",  "Acknowledge (Backchannel)")
    5   = @("sd", "Statement-non-opinion")
    10  = @("sv", "Statement-opinion")
    22  = @("%",  "Uninterpretable")
    25  = @("ba", "Appreciation")
    32  = @("qy", "Yes-No-Question")
    33  = @("sd", "Statement-non-opinion")
    46  = @("%",  "Uninterpretable")
    61  = @("sd", "Statement-non-opinion")
    69  = @("%",  "Uninterpretable")
    113 = @("aa", "Agree/Accept")
    114 = @("aa", "Agree/Accept")
    120 = @("ba", "Appreciation")
    124 = @("aa", "Agree/Accept")
    132 = @("ba", "Appreciation")
    133 = @("ba", "Appreciation")
    141 = @("aa", "Agree/Accept")
    147 = @("%",  "Uninterpretable")
    150 = @("sv", "Statement-opinion")
    165 = @("sd", "Statement-non-opinion")
    168 = @("sv", "Statement-opinion")
    181 = @("sd", "Statement-non-opinion")
    188 = @("sv", "Statement-opinion")
    213 = @("sd", "Statement-non-opinion")
    218 = @("sd", "Statement-non-opinion")
    222 = @("sd", "Statement-non-opinion")
    224 = @("sd", "Statement-non-opinion")
}

foreach ($row in $updates.Keys) {
    $values = $updates[$row]
    $ws.Cells.Item($row, 9).Value = $values[0]
    $ws.Cells.Item($row, 10).Value = $values[1]
}
